$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

$ws.Range("A2").Value = -79.3699
$ws.Range("B2").Value = -79.3344

$ws.Range("A3").Value = 33.1759
$ws.Range("B3").Value = 33.2054

$ws.Range("A4").Value = -79.0738
$ws.Range("B4").Value = -79.1095

$ws.Range("A5").Value = 33.4213
$ws.Range("B5").Value = 33.3918
